# Update "想去人数" (want-to-go count) values in column F on the
# "展览" and "全部类型" worksheets to reflect newly scraped data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1.xml) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value  = 2667
$wsExhibit.Range("F5").Value  = 1490
$wsExhibit.Range("F13").Value = 9087
$wsExhibit.Range("F18").Value = 179
$wsExhibit.Range("F25").Value = 2168
$wsExhibit.Range("F27").Value = 1875
$wsExhibit.Range("F31").Value = 1251
$wsExhibit.Range("F36").Value = 321
$wsExhibit.Range("F41").Value = 20
$wsExhibit.Range("F42").Value = 277
$wsExhibit.Range("F43").Value = 1370

# --- Sheet "全部类型" (sheet4.xml) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 2667
$wsAll.Range("F5").Value  = 1490
$wsAll.Range("F12").Value = 9087
$wsAll.Range("F19").Value = 179
$wsAll.Range("F24").Value = 2168
$wsAll.Range("F25").Value = 1875
$wsAll.Range("F28").Value = 1251
$wsAll.Range("F33").Value = 321
$wsAll.Range("F41").Value = 20
$wsAll.Range("F42").Value = 277
$wsAll.Range("F44").Value = 1370
